$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte planificación")

$ws.Range("G4").Value = 10
$ws.Range("J4").Value = 11
$ws.Range("G5").Value = 6
$ws.Range("J5").Value = 7
